# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.253.28"
$ws.Range("E2").Value = "  -0.51%  "

$ws.Range("D3").Value = "'1.804.56"
$ws.Range("E3").Value = "  -0.71%  "

$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").Value = "'314.26"
$ws.Range("E5").Value = "  -0.33%  "

$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").Value = "'0.5275"
$ws.Range("E7").Value = "  +3.15%  "

$ws.Range("D8").Value = "'0.3830"
$ws.Range("E8").Value = "  -2.92%  "

$ws.Range("D9").Value = "'0.08010"
$ws.Range("E9").Value = "  -0.58%  "

$ws.Range("D10").Value = "'41.44"
$ws.Range("E10").Value = "  -0.56%  "

$ws.Range("D11").Value = "'1.100"
$ws.Range("E11").Value = "  -0.56%  "

$ws.Range("D12").Value = "'6.325"
$ws.Range("E12").Value = "  +1.04%  "

$ws.Range("D13").Value = "'1.003"
$ws.Range("E13").Value = "  +0.15%  "

$ws.Range("D14").Value = "'20.60"
$ws.Range("E14").Value = "  -1.78%  "

$ws.Range("D15").Value = "'1.807.80"
$ws.Range("E15").Value = "  -0.97%  "

$ws.Range("D16").Value = "'7.325"
$ws.Range("E16").Value = "  -2.28%  "

$ws.Range("D17").Value = "'92.15"
$ws.Range("E17").Value = "  -0.53%  "

$ws.Range("D18").Value = "'0.00001096"
$ws.Range("E18").Value = "  -3.75%  "

$ws.Range("D19").Value = "'0.06612"
$ws.Range("E19").Value = "  -0.31%  "

$ws.Range("D20").Value = "'1.002"
$ws.Range("E20").Value = "  +0.09%  "

$ws.Range("D21").Value = "'17.38"
$ws.Range("E21").Value = "  -1.66%  "

$ws.Range("D22").Value = "'5.966"
$ws.Range("E22").Value = "  -1.96%  "

$ws.Range("D23").Value = "'28.313.89"
$ws.Range("E23").Value = "  -0.43%  "

$ws.Range("D24").Value = "'11.15"
$ws.Range("E24").Value = "  -0.88%  "

$ws.Range("D25").Value = "'2.256"
$ws.Range("E25").Value = "  -0.58%  "

$ws.Range("D26").Value = "'160.68"
$ws.Range("E26").Value = "  +3.71%  "

$ws.Range("D27").Value = "'20.47"
$ws.Range("E27").Value = "  -3.11%  "

$ws.Range("D28").Value = "'2.011.36"
$ws.Range("E28").Value = "  -1.06%  "

$ws.Range("D29").Value = "'2.359"
$ws.Range("E29").Value = "  -1.79%  "

$ws.Range("D30").Value = "'123.42"
$ws.Range("E30").Value = "  -1.97%  "

$ws.Range("D31").Value = "'0.1087"
$ws.Range("E31").Value = "  -1.27%  "

$ws.Range("E32").Value = "  -4.03%  "

$ws.Range("D33").Value = "'3.687"
$ws.Range("E33").Value = "  +0.97%  "

$ws.Range("D34").Value = "'5.558"
$ws.Range("E34").Value = "  -3.75%  "

$ws.Range("D35").Value = "'0.07256"
$ws.Range("E35").Value = "  +3.33%  "

$ws.Range("D36").Value = "'12.36"
$ws.Range("E36").Value = "  +9.52%  "

$ws.Range("E37").Value = "  -0.24%  "

$ws.Range("D38").Value = "'0.2159"
$ws.Range("E38").Value = "  -2.91%  "

$ws.Range("D39").Value = "'5.116"
$ws.Range("E39").Value = "  -1.72%  "

$ws.Range("D40").Value = "'8.665"
$ws.Range("E40").Value = "  -1.63%  "

$ws.Range("D41").Value = "'0.6207"
$ws.Range("E41").Value = "  -0.80%  "

$ws.Range("E42").Value = "  -0.69%  "

$ws.Range("D43").Value = "'1.370"
$ws.Range("E43").Value = "  -1.85%  "

$ws.Range("D44").Value = "'0.6022"
$ws.Range("E44").Value = "  +1.94%  "

$ws.Range("D45").Value = "'13.20"
$ws.Range("E45").Value = "  -2.18%  "

$ws.Range("D46").Value = "'3.767"
$ws.Range("E46").Value = "  +0.72%  "

$ws.Range("D47").Value = "'127.20"
$ws.Range("E47").Value = "  +1.88%  "

$ws.Range("D48").Value = "'1.217"
$ws.Range("E48").Value = "  +2.53%  "

$ws.Range("D49").Value = "'1.930"
$ws.Range("E49").Value = "  -2.18%  "

$ws.Range("D50").Value = "'0.06818"
$ws.Range("E50").Value = "  -0.98%  "

$ws.Range("D51").Value = "'73.30"
$ws.Range("E51").Value = "  -1.65%  "
